# Applies the text-run restructuring described in the commit
# "Test to see if change in powerpoint file triggers git changes".
# The visible text in every edited paragraph is unchanged; only the
# run (and, incidentally, endParaRPr) boundaries differ, consistent
# with someone merely clicking into placeholders and re-saving.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 10 ("GIT commands (2)- adding a file")
# Merge the 3 trailing runs of the "...ahead of origin/master by 1
# commit:" paragraph into a single run (keep the leading "w" run).
# ---------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(2)
$tr10 = $sh10.TextFrame.TextRange
$para10 = $tr10.Paragraphs(11, 1)
$c10 = $para10.Characters(2, $para10.Length - 2)
$c10.Text = "ill tell you that your local repository is ahead of origin/master by 1 commit:"

# ---------------------------------------------------------------
# Slide 11 ("GIT commands (3)- deleting a file")
# ---------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(2)
$tr11 = $sh11.TextFrame.TextRange

# Paragraph 1: "Now to remove the test.txt " + "file" -> one run
$para11a = $tr11.Paragraphs(1, 1)
$c11a = $para11a.Characters(1, $para11a.Length)
$c11a.Text = "Now to remove the test.txt file"

# Paragraph 11: "...should be gone" + "." -> one run
$para11b = $tr11.Paragraphs(11, 1)
$c11b = $para11b.Characters(1, $para11b.Length)
$c11b.Text = "Will update the remote repository on GitHub. Look at the webpage to confirm this. The file test.txt should be gone."

# ---------------------------------------------------------------
# Slide 2 ("What is version control?")
# ---------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(4)
$tr2 = $sh2.TextFrame.TextRange

# Paragraph 1: split "making copies " out of the run so it reads as
# "If you are " | "making " | "copies " | "of projects and renaming..."
$para2a = $tr2.Paragraphs(1, 1)
$cDel = $para2a.Characters(12, 14)
$cDel.Delete()
$afterIf = $para2a.Characters(1, 11)
$insMaking = $afterIf.InsertAfter("making ")
$insCopies = $insMaking.InsertAfter("copies ")

# Paragraph 5: "If you want to try something new (no matter how " +
# "silly" + ") " + "with your code but don't" -> one run
$para2b = $tr2.Paragraphs(5, 1)
$c2b = $para2b.Characters(1, 74)
$c2b.Text = "If you want to try something new (no matter how silly) with your code but don’t"

# ---------------------------------------------------------------
# Slide 5 ("GIT first steps")
# Merge " " + "yours (Pull request in GitHub)." into a single run.
# ---------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(2)
$tr5 = $sh5.TextFrame.TextRange
$para5 = $tr5.Paragraphs(4, 1)
$c5 = $para5.Characters($para5.Length - 32, 33)
$c5.Text = " yours (Pull request in GitHub)."
